# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.981.28'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '1.897.61'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '''0.8319'
$ws.Range("E5").Value = '  +4.61%  '
$ws.Range("D6").Value = '''241.90'
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '''0.3284'
$ws.Range("E8").Value = '  +3.01%  '
$ws.Range("D9").Value = '''26.57'
$ws.Range("E9").Value = '  +0.63%  '
$ws.Range("D10").Value = '''0.07041'
$ws.Range("E10").Value = '  +1.01%  '
$ws.Range("D11").Value = '''0.08080'
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("D12").Value = '''0.7606'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").Value = '1.897.87'
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").Value = '''5.245'
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("D15").Value = '''92.19'
$ws.Range("E15").Value = '  -1.69%  '
$ws.Range("D16").Value = '29.987.66'
$ws.Range("E16").Value = '  -0.57%  '
$ws.Range("D17").Value = '''14.12'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").Value = '''5.865'
$ws.Range("E18").Value = '  -2.44%  '
$ws.Range("D19").Value = '''244.24'
$ws.Range("E19").Value = '  -2.16%  '
$ws.Range("D20").Value = '''0.000007752'
$ws.Range("E20").Value = '  -1.10%  '
$ws.Range("D21").Value = '''1.000'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = '2.151.99'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").Value = '''6.960'
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").Value = '''0.1744'
$ws.Range("E25").Value = '  +25.04%  '
$ws.Range("D26").Value = '''9.252'
$ws.Range("E26").Value = '  -1.06%  '
$ws.Range("D27").Value = '''165.63'
$ws.Range("E27").Value = '  -2.08%  '
$ws.Range("D28").Value = '''18.89'
$ws.Range("E28").Value = '  -0.70%  '
$ws.Range("D29").Value = '''2.092'
$ws.Range("E29").Value = '  +1.57%  '
$ws.Range("D30").Value = '''1.360'
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("D31").Value = '''1.516'
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").Value = '''0.05954'
$ws.Range("E32").Value = '  +10.60%  '
$ws.Range("D33").Value = '''4.279'
$ws.Range("E33").Value = '  -2.22%  '
$ws.Range("D34").Value = '''4.070'
$ws.Range("E34").Value = '  -1.42%  '
$ws.Range("D35").Value = '''1.268'
$ws.Range("E35").Value = '  -0.33%  '
$ws.Range("D36").Value = '''0.7310'
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("D37").Value = '''2.722'
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D39").Value = '''2.778'
$ws.Range("D40").Value = '''0.4439'
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("D41").Value = '''72.39'
$ws.Range("E41").Value = '  -0.59%  '
$ws.Range("D42").Value = '''5.850'
$ws.Range("E42").Value = '  -5.59%  '
$ws.Range("D43").Value = '''0.8511'
$ws.Range("E43").Value = '  +1.85%  '
$ws.Range("D44").Value = '''1.000'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '''1.900'
$ws.Range("E45").Value = '  -0.34%  '
$ws.Range("D46").Value = '''101.89'
$ws.Range("E46").Value = '  +0.99%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''9.824'
$ws.Range("E47").Value = '  -0.93%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '''7.547'
$ws.Range("E48").Value = '  -1.18%  '
$ws.Range("D49").Value = '''988.04'
$ws.Range("E49").Value = '  +2.33%  '
$ws.Range("D50").Value = '2.044.67'
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''1.519'
$ws.Range("E51").Value = '  +0.08%  '
